$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.004.07'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '2.581.01'
$ws.Range('E3').Value = '  -3.78%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '550.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.88'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +1.85%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '3.037.88'
$ws.Range('E13').Value = '  -3.86%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.62'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '61.915.54'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '2.585.64'
$ws.Range('E17').Value = '  -3.86%  '
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '338.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').Value = '  -4.13%  '
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.168'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.22'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.28'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.10%  '
$ws.Range('D29').Value = '0.0₃0839'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E30').Value = '  +2.83%  '
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '162.92'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '19.23'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '329.91'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.06'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.907'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.76%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.95'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '37.64'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.98'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.24%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').Value = '2.112.11'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '19.56'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0967'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0239'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.67%  '
